$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Nodo" key-exchange block (rows 13-14, columns O-R)
$ws.Range("O13").Value = "Nodo"
$ws.Range("P13").Value = 6
$ws.Range("Q13").Value = 2
$ws.Range("R13").Value = 32

$ws.Range("P14").Value = "MAC"
$ws.Range("Q14").Value = "NodeId"
$ws.Range("R14").Value = "Key"

# Zoom out the view and move the selection, as in the author's edit
$excel.ActiveWindow.Zoom = 40
[void]$ws.Range("T18").Select()
